$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (e.g. AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Re-set the values since PasteSpecial formats only shouldn't touch them, but just in case
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in team record data for each player row (2 through 54)
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = 85   # AD = Wins
    $ws.Cells.Item($r, 31).Value = 77   # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = Ties
}
